$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.223433242506812
$ws.Range("C2").Value = 0.4904632152588556
$ws.Range("J2").Value = 0.01362397820163488
$ws.Range("P2").Value = 0.1771117166212534
$ws.Range("S2").Value = 0.09536784741144415
$ws.Range("C3").Value = 0.00546448087431694
$ws.Range("J3").Value = 0.06557377049180328
$ws.Range("P3").Value = 0.7595628415300546
$ws.Range("S3").Value = 0.1693989071038251
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.6444444444444445
$ws.Range("S4").Value = 0.3111111111111111
$ws.Range("B6").Value = 0.06532663316582915
$ws.Range("F6").Value = 0.06030150753768844
$ws.Range("J6").Value = 0.3417085427135678
$ws.Range("O6").Value = 0.02010050251256281
$ws.Range("Q6").Value = 0.1055276381909548
$ws.Range("R6").Value = 0.03517587939698492
$ws.Range("S6").Value = 0.3718592964824121
$ws.Range("B7").Value = 0.1473684210526316
$ws.Range("D7").Value = 0.02105263157894737
$ws.Range("F7").Value = 0.03157894736842105
$ws.Range("J7").Value = 0.1736842105263158
$ws.Range("O7").Value = 0.03157894736842105
$ws.Range("Q7").Value = 0.1263157894736842
$ws.Range("R7").Value = 0.05789473684210526
$ws.Range("S7").Value = 0.4105263157894737
$ws.Range("B8").Value = 0.1156862745098039
$ws.Range("D8").Value = 0.02156862745098039
$ws.Range("F8").Value = 0.04313725490196078
$ws.Range("J8").Value = 0.1411764705882353
$ws.Range("O8").Value = 0.01568627450980392
$ws.Range("Q8").Value = 0.1843137254901961
$ws.Range("R8").Value = 0.09215686274509804
$ws.Range("S8").Value = 0.3862745098039216
$ws.Range("B9").Value = 0.1260504201680672
$ws.Range("D9").Value = 0.03361344537815126
$ws.Range("F9").Value = 0.07563025210084033
$ws.Range("J9").Value = 0.08403361344537816
$ws.Range("O9").Value = 0.04201680672268908
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.03361344537815126
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.1228464419475655
$ws.Range("D10").Value = 0.02022471910112359
$ws.Range("E10").Value = 0.000749063670411985
$ws.Range("F10").Value = 0.06142322097378277
$ws.Range("J10").Value = 0.1415730337078652
$ws.Range("O10").Value = 0.01722846441947565
$ws.Range("Q10").Value = 0.2411985018726592
$ws.Range("R10").Value = 0.07415730337078652
$ws.Range("S10").Value = 0.3205992509363296
$ws.Range("G11").Value = 0.1893939393939394
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.196969696969697
$ws.Range("L11").Value = 0.5303030303030303
$ws.Range("G12").Value = 0.7837837837837838
$ws.Range("J12").Value = 0.1689189189189189
$ws.Range("L12").Value = 0.03378378378378379
$ws.Range("S12").Value = 0.01351351351351351
$ws.Range("G13").Value = 0.7027027027027027
$ws.Range("J13").Value = 0.2702702702702703
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("F15").Value = 0.02061855670103093
$ws.Range("H15").Value = 0.154639175257732
$ws.Range("I15").Value = 0.05154639175257732
$ws.Range("J15").Value = 0.3969072164948453
$ws.Range("K15").Value = 0.04639175257731959
$ws.Range("O15").Value = 0.03608247422680412
$ws.Range("S15").Value = 0.2938144329896907
$ws.Range("F16").Value = 0.01282051282051282
$ws.Range("H16").Value = 0.2008547008547009
$ws.Range("I16").Value = 0.05555555555555555
$ws.Range("J16").Value = 0.452991452991453
$ws.Range("K16").Value = 0.09829059829059829
$ws.Range("M16").Value = 0.008547008547008548
$ws.Range("O16").Value = 0.04273504273504274
$ws.Range("S16").Value = 0.1282051282051282
$ws.Range("F17").Value = 0.02114164904862579
$ws.Range("H17").Value = 0.1902748414376321
$ws.Range("I17").Value = 0.05708245243128964
$ws.Range("J17").Value = 0.4778012684989429
$ws.Range("K17").Value = 0.08456659619450317
$ws.Range("M17").Value = 0.008456659619450317
$ws.Range("O17").Value = 0.040169133192389
$ws.Range("S17").Value = 0.120507399577167
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.1987951807228916
$ws.Range("I18").Value = 0.06626506024096386
$ws.Range("J18").Value = 0.463855421686747
$ws.Range("K18").Value = 0.07228915662650602
$ws.Range("M18").Value = 0.006024096385542169
$ws.Range("O18").Value = 0.04819277108433735
$ws.Range("S18").Value = 0.1325301204819277
$ws.Range("F19").Value = 0.01978565539983512
$ws.Range("H19").Value = 0.2629843363561418
$ws.Range("I19").Value = 0.0494641384995878
$ws.Range("J19").Value = 0.3446001648804617
$ws.Range("K19").Value = 0.1030502885408079
$ws.Range("M19").Value = 0.02555647155812036
$ws.Range("O19").Value = 0.06924979389942292
$ws.Range("S19").Value = 0.1253091508656224
